$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$dateFmt = $ws.Cells.Item(700, 4).NumberFormat

# Row 749
$ws.Cells.Item(749, 1).Value = 11
$ws.Cells.Item(749, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(749, 3).Value = 'Bíobío'
$ws.Cells.Item(749, 4).Value = 44628
$ws.Cells.Item(749, 4).NumberFormat = $dateFmt
$ws.Cells.Item(749, 5).Value = 8
$ws.Cells.Item(749, 6).Value = 'Fruta'
$ws.Cells.Item(749, 7).Value = 100104
$ws.Cells.Item(749, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(749, 9).Value = 100104002
$ws.Cells.Item(749, 10).Value = 'Manzana'
$ws.Cells.Item(749, 11).Value = 'Granny Smith'
$ws.Cells.Item(749, 12).Value = 'Especial'
$ws.Cells.Item(749, 13).Value = 50
$ws.Cells.Item(749, 14).Value = 12000
$ws.Cells.Item(749, 15).Value = 12000
$ws.Cells.Item(749, 16).Value = 12000
$ws.Cells.Item(749, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(749, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(749, 19).Value = 750
$ws.Cells.Item(749, 20).Value = 16

# Row 750
$ws.Cells.Item(750, 1).Value = 11
$ws.Cells.Item(750, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(750, 3).Value = 'Bíobío'
$ws.Cells.Item(750, 4).Value = 44628
$ws.Cells.Item(750, 4).NumberFormat = $dateFmt
$ws.Cells.Item(750, 5).Value = 8
$ws.Cells.Item(750, 6).Value = 'Fruta'
$ws.Cells.Item(750, 7).Value = 100104
$ws.Cells.Item(750, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(750, 9).Value = 100104002
$ws.Cells.Item(750, 10).Value = 'Manzana'
$ws.Cells.Item(750, 11).Value = 'Granny Smith'
$ws.Cells.Item(750, 12).Value = 'Primera'
$ws.Cells.Item(750, 13).Value = 100
$ws.Cells.Item(750, 14).Value = 10000
$ws.Cells.Item(750, 15).Value = 10000
$ws.Cells.Item(750, 16).Value = 10000
$ws.Cells.Item(750, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(750, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(750, 19).Value = 625
$ws.Cells.Item(750, 20).Value = 16

# Row 751
$ws.Cells.Item(751, 1).Value = 11
$ws.Cells.Item(751, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(751, 3).Value = 'Bíobío'
$ws.Cells.Item(751, 4).Value = 44628
$ws.Cells.Item(751, 4).NumberFormat = $dateFmt
$ws.Cells.Item(751, 5).Value = 8
$ws.Cells.Item(751, 6).Value = 'Fruta'
$ws.Cells.Item(751, 7).Value = 100104
$ws.Cells.Item(751, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(751, 9).Value = 100104002
$ws.Cells.Item(751, 10).Value = 'Manzana'
$ws.Cells.Item(751, 11).Value = 'Granny Smith'
$ws.Cells.Item(751, 12).Value = 'Segunda'
$ws.Cells.Item(751, 13).Value = 100
$ws.Cells.Item(751, 14).Value = 8000
$ws.Cells.Item(751, 15).Value = 8000
$ws.Cells.Item(751, 16).Value = 8000
$ws.Cells.Item(751, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(751, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(751, 19).Value = 500
$ws.Cells.Item(751, 20).Value = 16

# Row 752
$ws.Cells.Item(752, 1).Value = 11
$ws.Cells.Item(752, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(752, 3).Value = 'Bíobío'
$ws.Cells.Item(752, 4).Value = 44628
$ws.Cells.Item(752, 4).NumberFormat = $dateFmt
$ws.Cells.Item(752, 5).Value = 8
$ws.Cells.Item(752, 6).Value = 'Fruta'
$ws.Cells.Item(752, 7).Value = 100104
$ws.Cells.Item(752, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(752, 9).Value = 100104002
$ws.Cells.Item(752, 10).Value = 'Manzana'
$ws.Cells.Item(752, 11).Value = 'Royal Gala'
$ws.Cells.Item(752, 12).Value = 'Especial'
$ws.Cells.Item(752, 13).Value = 50
$ws.Cells.Item(752, 14).Value = 12000
$ws.Cells.Item(752, 15).Value = 12000
$ws.Cells.Item(752, 16).Value = 12000
$ws.Cells.Item(752, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(752, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(752, 19).Value = 750
$ws.Cells.Item(752, 20).Value = 16

# Row 753
$ws.Cells.Item(753, 1).Value = 11
$ws.Cells.Item(753, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(753, 3).Value = 'Bíobío'
$ws.Cells.Item(753, 4).Value = 44628
$ws.Cells.Item(753, 4).NumberFormat = $dateFmt
$ws.Cells.Item(753, 5).Value = 8
$ws.Cells.Item(753, 6).Value = 'Fruta'
$ws.Cells.Item(753, 7).Value = 100104
$ws.Cells.Item(753, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(753, 9).Value = 100104002
$ws.Cells.Item(753, 10).Value = 'Manzana'
$ws.Cells.Item(753, 11).Value = 'Royal Gala'
$ws.Cells.Item(753, 12).Value = 'Primera'
$ws.Cells.Item(753, 13).Value = 100
$ws.Cells.Item(753, 14).Value = 10000
$ws.Cells.Item(753, 15).Value = 10000
$ws.Cells.Item(753, 16).Value = 10000
$ws.Cells.Item(753, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(753, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(753, 19).Value = 625
$ws.Cells.Item(753, 20).Value = 16

# Row 754
$ws.Cells.Item(754, 1).Value = 11
$ws.Cells.Item(754, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(754, 3).Value = 'Bíobío'
$ws.Cells.Item(754, 4).Value = 44628
$ws.Cells.Item(754, 4).NumberFormat = $dateFmt
$ws.Cells.Item(754, 5).Value = 8
$ws.Cells.Item(754, 6).Value = 'Fruta'
$ws.Cells.Item(754, 7).Value = 100104
$ws.Cells.Item(754, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(754, 9).Value = 100104002
$ws.Cells.Item(754, 10).Value = 'Manzana'
$ws.Cells.Item(754, 11).Value = 'Royal Gala'
$ws.Cells.Item(754, 12).Value = 'Segunda'
$ws.Cells.Item(754, 13).Value = 100
$ws.Cells.Item(754, 14).Value = 8000
$ws.Cells.Item(754, 15).Value = 8000
$ws.Cells.Item(754, 16).Value = 8000
$ws.Cells.Item(754, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(754, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(754, 19).Value = 500
$ws.Cells.Item(754, 20).Value = 16

# Row 755
$ws.Cells.Item(755, 1).Value = 11
$ws.Cells.Item(755, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(755, 3).Value = 'Bíobío'
$ws.Cells.Item(755, 4).Value = 44335
$ws.Cells.Item(755, 4).NumberFormat = $dateFmt
$ws.Cells.Item(755, 5).Value = 8
$ws.Cells.Item(755, 6).Value = 'Fruta'
$ws.Cells.Item(755, 7).Value = 100104
$ws.Cells.Item(755, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(755, 9).Value = 100104002
$ws.Cells.Item(755, 10).Value = 'Manzana'
$ws.Cells.Item(755, 11).Value = 'Fuji royal'
$ws.Cells.Item(755, 12).Value = 'Especial'
$ws.Cells.Item(755, 13).Value = 50
$ws.Cells.Item(755, 14).Value = 12000
$ws.Cells.Item(755, 15).Value = 12000
$ws.Cells.Item(755, 16).Value = 12000
$ws.Cells.Item(755, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(755, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(755, 19).Value = 750
$ws.Cells.Item(755, 20).Value = 16

# Row 756
$ws.Cells.Item(756, 1).Value = 11
$ws.Cells.Item(756, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(756, 3).Value = 'Bíobío'
$ws.Cells.Item(756, 4).Value = 44335
$ws.Cells.Item(756, 4).NumberFormat = $dateFmt
$ws.Cells.Item(756, 5).Value = 8
$ws.Cells.Item(756, 6).Value = 'Fruta'
$ws.Cells.Item(756, 7).Value = 100104
$ws.Cells.Item(756, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(756, 9).Value = 100104002
$ws.Cells.Item(756, 10).Value = 'Manzana'
$ws.Cells.Item(756, 11).Value = 'Fuji royal'
$ws.Cells.Item(756, 12).Value = 'Primera'
$ws.Cells.Item(756, 13).Value = 50
$ws.Cells.Item(756, 14).Value = 10000
$ws.Cells.Item(756, 15).Value = 10000
$ws.Cells.Item(756, 16).Value = 10000
$ws.Cells.Item(756, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(756, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(756, 19).Value = 625
$ws.Cells.Item(756, 20).Value = 16

# Row 757
$ws.Cells.Item(757, 1).Value = 11
$ws.Cells.Item(757, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(757, 3).Value = 'Bíobío'
$ws.Cells.Item(757, 4).Value = 44335
$ws.Cells.Item(757, 4).NumberFormat = $dateFmt
$ws.Cells.Item(757, 5).Value = 8
$ws.Cells.Item(757, 6).Value = 'Fruta'
$ws.Cells.Item(757, 7).Value = 100104
$ws.Cells.Item(757, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(757, 9).Value = 100104002
$ws.Cells.Item(757, 10).Value = 'Manzana'
$ws.Cells.Item(757, 11).Value = 'Fuji royal'
$ws.Cells.Item(757, 12).Value = 'Segunda'
$ws.Cells.Item(757, 13).Value = 50
$ws.Cells.Item(757, 14).Value = 8000
$ws.Cells.Item(757, 15).Value = 8000
$ws.Cells.Item(757, 16).Value = 8000
$ws.Cells.Item(757, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(757, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(757, 19).Value = 500
$ws.Cells.Item(757, 20).Value = 16

# Row 758
$ws.Cells.Item(758, 1).Value = 11
$ws.Cells.Item(758, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(758, 3).Value = 'Bíobío'
$ws.Cells.Item(758, 4).Value = 44335
$ws.Cells.Item(758, 4).NumberFormat = $dateFmt
$ws.Cells.Item(758, 5).Value = 8
$ws.Cells.Item(758, 6).Value = 'Fruta'
$ws.Cells.Item(758, 7).Value = 100104
$ws.Cells.Item(758, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(758, 9).Value = 100104002
$ws.Cells.Item(758, 10).Value = 'Manzana'
$ws.Cells.Item(758, 11).Value = 'Granny Smith'
$ws.Cells.Item(758, 12).Value = 'Especial'
$ws.Cells.Item(758, 13).Value = 50
$ws.Cells.Item(758, 14).Value = 12000
$ws.Cells.Item(758, 15).Value = 12000
$ws.Cells.Item(758, 16).Value = 12000
$ws.Cells.Item(758, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(758, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(758, 19).Value = 750
$ws.Cells.Item(758, 20).Value = 16

# Row 759
$ws.Cells.Item(759, 1).Value = 11
$ws.Cells.Item(759, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(759, 3).Value = 'Bíobío'
$ws.Cells.Item(759, 4).Value = 44335
$ws.Cells.Item(759, 4).NumberFormat = $dateFmt
$ws.Cells.Item(759, 5).Value = 8
$ws.Cells.Item(759, 6).Value = 'Fruta'
$ws.Cells.Item(759, 7).Value = 100104
$ws.Cells.Item(759, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(759, 9).Value = 100104002
$ws.Cells.Item(759, 10).Value = 'Manzana'
$ws.Cells.Item(759, 11).Value = 'Granny Smith'
$ws.Cells.Item(759, 12).Value = 'Primera'
$ws.Cells.Item(759, 13).Value = 50
$ws.Cells.Item(759, 14).Value = 10000
$ws.Cells.Item(759, 15).Value = 10000
$ws.Cells.Item(759, 16).Value = 10000
$ws.Cells.Item(759, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(759, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(759, 19).Value = 625
$ws.Cells.Item(759, 20).Value = 16

# Row 760
$ws.Cells.Item(760, 1).Value = 11
$ws.Cells.Item(760, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(760, 3).Value = 'Bíobío'
$ws.Cells.Item(760, 4).Value = 44335
$ws.Cells.Item(760, 4).NumberFormat = $dateFmt
$ws.Cells.Item(760, 5).Value = 8
$ws.Cells.Item(760, 6).Value = 'Fruta'
$ws.Cells.Item(760, 7).Value = 100104
$ws.Cells.Item(760, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(760, 9).Value = 100104002
$ws.Cells.Item(760, 10).Value = 'Manzana'
$ws.Cells.Item(760, 11).Value = 'Granny Smith'
$ws.Cells.Item(760, 12).Value = 'Segunda'
$ws.Cells.Item(760, 13).Value = 50
$ws.Cells.Item(760, 14).Value = 8000
$ws.Cells.Item(760, 15).Value = 8000
$ws.Cells.Item(760, 16).Value = 8000
$ws.Cells.Item(760, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(760, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(760, 19).Value = 500
$ws.Cells.Item(760, 20).Value = 16

# Row 761
$ws.Cells.Item(761, 1).Value = 11
$ws.Cells.Item(761, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(761, 3).Value = 'Bíobío'
$ws.Cells.Item(761, 4).Value = 44160
$ws.Cells.Item(761, 4).NumberFormat = $dateFmt
$ws.Cells.Item(761, 5).Value = 8
$ws.Cells.Item(761, 6).Value = 'Fruta'
$ws.Cells.Item(761, 7).Value = 100104
$ws.Cells.Item(761, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(761, 9).Value = 100104002
$ws.Cells.Item(761, 10).Value = 'Manzana'
$ws.Cells.Item(761, 11).Value = 'Granny Smith'
$ws.Cells.Item(761, 12).Value = 'Primera'
$ws.Cells.Item(761, 13).Value = 200
$ws.Cells.Item(761, 14).Value = 10000
$ws.Cells.Item(761, 15).Value = 11000
$ws.Cells.Item(761, 16).Value = 10500
$ws.Cells.Item(761, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(761, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(761, 19).Value = 656
$ws.Cells.Item(761, 20).Value = 16

# Row 762
$ws.Cells.Item(762, 1).Value = 11
$ws.Cells.Item(762, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(762, 3).Value = 'Bíobío'
$ws.Cells.Item(762, 4).Value = 44160
$ws.Cells.Item(762, 4).NumberFormat = $dateFmt
$ws.Cells.Item(762, 5).Value = 8
$ws.Cells.Item(762, 6).Value = 'Fruta'
$ws.Cells.Item(762, 7).Value = 100104
$ws.Cells.Item(762, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(762, 9).Value = 100104002
$ws.Cells.Item(762, 10).Value = 'Manzana'
$ws.Cells.Item(762, 11).Value = 'Granny Smith'
$ws.Cells.Item(762, 12).Value = 'Segunda'
$ws.Cells.Item(762, 13).Value = 100
$ws.Cells.Item(762, 14).Value = 9000
$ws.Cells.Item(762, 15).Value = 9000
$ws.Cells.Item(762, 16).Value = 9000
$ws.Cells.Item(762, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(762, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(762, 19).Value = 562
$ws.Cells.Item(762, 20).Value = 16

# Row 763
$ws.Cells.Item(763, 1).Value = 11
$ws.Cells.Item(763, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(763, 3).Value = 'Bíobío'
$ws.Cells.Item(763, 4).Value = 44160
$ws.Cells.Item(763, 4).NumberFormat = $dateFmt
$ws.Cells.Item(763, 5).Value = 8
$ws.Cells.Item(763, 6).Value = 'Fruta'
$ws.Cells.Item(763, 7).Value = 100104
$ws.Cells.Item(763, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(763, 9).Value = 100104002
$ws.Cells.Item(763, 10).Value = 'Manzana'
$ws.Cells.Item(763, 11).Value = 'Pink Lady'
$ws.Cells.Item(763, 12).Value = 'Primera'
$ws.Cells.Item(763, 13).Value = 200
$ws.Cells.Item(763, 14).Value = 10000
$ws.Cells.Item(763, 15).Value = 11000
$ws.Cells.Item(763, 16).Value = 10500
$ws.Cells.Item(763, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(763, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(763, 19).Value = 656
$ws.Cells.Item(763, 20).Value = 16

# Row 764
$ws.Cells.Item(764, 1).Value = 11
$ws.Cells.Item(764, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(764, 3).Value = 'Bíobío'
$ws.Cells.Item(764, 4).Value = 44160
$ws.Cells.Item(764, 4).NumberFormat = $dateFmt
$ws.Cells.Item(764, 5).Value = 8
$ws.Cells.Item(764, 6).Value = 'Fruta'
$ws.Cells.Item(764, 7).Value = 100104
$ws.Cells.Item(764, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(764, 9).Value = 100104002
$ws.Cells.Item(764, 10).Value = 'Manzana'
$ws.Cells.Item(764, 11).Value = 'Pink Lady'
$ws.Cells.Item(764, 12).Value = 'Segunda'
$ws.Cells.Item(764, 13).Value = 100
$ws.Cells.Item(764, 14).Value = 9000
$ws.Cells.Item(764, 15).Value = 9000
$ws.Cells.Item(764, 16).Value = 9000
$ws.Cells.Item(764, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(764, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(764, 19).Value = 562
$ws.Cells.Item(764, 20).Value = 16
